# Auto-generated edit script: updates crypto price/volume columns
# and swaps the ordi/Algorand rows (45/46) per the commit diff.
#
# D-column (Price) values are prefixed with a leading apostrophe so
# Excel stores them as Text (matching the original t="inlineStr" cells)
# instead of auto-converting number-looking strings (e.g. "100.10")
# into numeric values that would drop significant trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.952.74"
$ws.Range("E2").Value = "  +1.56%  "

$ws.Range("D3").Value = "'2.243.33"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("D5").Value = "'319.07"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").Value = "'100.10"
$ws.Range("E6").Value = "  +1.43%  "

$ws.Range("D7").Value = "'0.573"
$ws.Range("E7").Value = "  -1.32%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  -2.92%  "

$ws.Range("D10").Value = "'36.60"
$ws.Range("E10").Value = "  -1.27%  "

$ws.Range("D11").Value = "'0.0824"
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("D12").Value = "'7.48"
$ws.Range("E12").Value = "  -2.77%  "

$ws.Range("E13").Value = "  -1.88%  "

$ws.Range("D14").Value = "'2.585.05"
$ws.Range("E14").Value = "  +0.72%  "

$ws.Range("D15").Value = "'0.847"
$ws.Range("E15").Value = "  -1.76%  "

$ws.Range("D16").Value = "'14.24"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").Value = "'2.248.04"
$ws.Range("E17").Value = "  +0.73%  "

$ws.Range("D18").Value = "'43.864.54"
$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("D19").Value = "'13.28"
$ws.Range("E19").Value = "  -5.14%  "

$ws.Range("D20").Value = "'0.0₃0971"
$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("D21").Value = "'6.40"
$ws.Range("E21").Value = "  -2.29%  "

$ws.Range("D22").Value = "'65.17"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "'3.08"
$ws.Range("E23").Value = "  -4.05%  "

$ws.Range("D24").Value = "'233.13"
$ws.Range("E24").Value = "  -1.55%  "

$ws.Range("E25").Value = "  -5.73%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").Value = "'10.55"
$ws.Range("E27").Value = "  +5.11%  "

$ws.Range("D28").Value = "'38.31"
$ws.Range("E28").Value = "  +4.64%  "

$ws.Range("E29").Value = "  -1.40%  "

$ws.Range("D30").Value = "'6.04"
$ws.Range("E30").Value = "  -5.71%  "

$ws.Range("D31").Value = "'158.06"
$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("D32").Value = "'20.02"
$ws.Range("E32").Value = "  -1.09%  "

$ws.Range("E33").Value = "  -3.04%  "

$ws.Range("E34").Value = "  +0.37%  "

$ws.Range("E35").Value = "  -4.97%  "

$ws.Range("D36").Value = "'0.112"
$ws.Range("E36").Value = "  +7.07%  "

$ws.Range("E37").Value = "  +5.18%  "

$ws.Range("E38").Value = "  -2.22%  "

$ws.Range("D39").Value = "'16.15"
$ws.Range("E39").Value = "  +11.98%  "

$ws.Range("D40").Value = "'3.65"
$ws.Range("E40").Value = "  -1.69%  "

$ws.Range("D41").Value = "'4.13"
$ws.Range("E41").Value = "  -6.13%  "

$ws.Range("D42").Value = "'0.0312"
$ws.Range("E42").Value = "  -2.14%  "

$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").Value = "'1.759.13"
$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.194"
$ws.Range("E45").Value = "  -4.15%  "

$ws.Range("B46").Value = "ordi"
$ws.Range("C46").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D46").Value = "'73.92"
$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("D47").Value = "'80.79"
$ws.Range("E47").Value = "  -3.55%  "

$ws.Range("D48").Value = "'5.14"
$ws.Range("E48").Value = "  -2.66%  "

$ws.Range("D49").Value = "'102.91"
$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("E50").Value = "  +0.95%  "

$ws.Range("D51").Value = "'56.96"
$ws.Range("E51").Value = "  -1.96%  "
